$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 857.8182
$ws.Range("I34").Value = 857.8182
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 857.8182
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -654.8182
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 857.8182
$ws.Range("I36").Value = 857.8182
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 857.8182
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -142.8182
$ws.Range("N36").ClearContents()
$ws.Range("H94").Value = 10876.733
$ws.Range("I94").Value = 10876.733
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 10876.733
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -10425.733
$ws.Range("N94").ClearContents()
$ws.Range("H106").Value = 1840.7273
$ws.Range("J106").Value = 4500
$ws.Range("L106").Value = 4500
$ws.Range("N106").Value = -5762
$ws.Range("H133").Value = 200000
$ws.Range("J133").Value = 200000
$ws.Range("L133").Value = 200000
$ws.Range("N133").Value = -210120
$ws.Range("H137").Value = 3935.25
$ws.Range("I137").Value = 1914
$ws.Range("J137").Value = 9999
$ws.Range("K137").Value = 5742
$ws.Range("L137").Value = 29997
$ws.Range("M137").Value = -3192
$ws.Range("N137").Value = -35097
$ws.Range("H138").Value = 4425.2354
$ws.Range("J138").Value = 4710.769
$ws.Range("L138").Value = 14132.307
$ws.Range("N138").Value = -24412.307

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6271492.5
$ws.Range("I32").Value = 15688.1
$ws.Range("J32").Value = 16697833
$ws.Range("K32").Value = 15688.1
$ws.Range("L32").Value = 16697833
$ws.Range("M32").Value = -15401.1
$ws.Range("N32").Value = -16698407
$ws.Range("H61").Value = 5877.9165
$ws.Range("I61").Value = 5392.778
$ws.Range("J61").Value = 7333.3335
$ws.Range("K61").Value = 5392.778
$ws.Range("L61").Value = 7333.3335
$ws.Range("M61").Value = -5180.778
$ws.Range("N61").Value = -7757.3335
$ws.Range("H63").Value = 3462.5557
$ws.Range("J63").Value = 10000
$ws.Range("L63").Value = 10000
$ws.Range("N63").Value = -11372
$ws.Range("H66").Value = 3462.5557
$ws.Range("J66").Value = 10000
$ws.Range("L66").Value = 50000
$ws.Range("N66").Value = -56864
$ws.Range("H96").Value = 14111.125
$ws.Range("J96").Value = 14111.125
$ws.Range("L96").Value = 14111.125
$ws.Range("N96").Value = -19603.125
$ws.Range("H136").Value = 5877.9165
$ws.Range("I136").Value = 5392.778
$ws.Range("J136").Value = 7333.3335
$ws.Range("K136").Value = 16178.334
$ws.Range("L136").Value = 22000.0005
$ws.Range("M136").Value = -13628.334
$ws.Range("N136").Value = -27100.0005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8550.5
$ws.Range("I20").Value = 7949.3335
$ws.Range("K20").Value = 7949.3335
$ws.Range("M20").Value = -7702.3335
$ws.Range("H86").Value = 5916.25
$ws.Range("I86").Value = 2999.4
$ws.Range("J86").Value = 7999.7144
$ws.Range("K86").Value = 2999.4
$ws.Range("L86").Value = 7999.7144
$ws.Range("M86").Value = -1876.4
$ws.Range("N86").Value = -10245.7144
$ws.Range("H89").Value = 5916.25
$ws.Range("I89").Value = 2999.4
$ws.Range("J89").Value = 7999.7144
$ws.Range("K89").Value = 14997
$ws.Range("L89").Value = 39998.572
$ws.Range("M89").Value = -9381
$ws.Range("N89").Value = -51230.572
$ws.Range("H97").Value = 33749.5
$ws.Range("I97").Value = 24999.334
$ws.Range("K97").Value = 24999.334
$ws.Range("M97").Value = -24008.334
$ws.Range("H139").Value = 70000
$ws.Range("I139").Value = 70000
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 70000
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -64860
$ws.Range("N139").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 15547.5
$ws.Range("I103").Value = 15547.5
$ws.Range("K103").Value = 15547.5
$ws.Range("M103").Value = -14375.5
$ws.Range("H104").Value = 56333.332
$ws.Range("I104").Value = 45000
$ws.Range("J104").Value = 62000
$ws.Range("K104").Value = 45000
$ws.Range("L104").Value = 62000
$ws.Range("M104").Value = -42379
$ws.Range("N104").Value = -67242
$ws.Range("H141").Value = 49867.91
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 49867.91
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 49867.91
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -60227.91

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1441.25
$ws.Range("I18").Value = 930.1667
$ws.Range("K18").Value = 2790.5001
$ws.Range("M18").Value = -2621.5001
$ws.Range("H36").Value = 400
$ws.Range("I36").Value = 350
$ws.Range("J36").Value = 500
$ws.Range("K36").Value = 1050
$ws.Range("L36").Value = 1500
$ws.Range("M36").Value = -881
$ws.Range("N36").Value = -1838
$ws.Range("H38").Value = 244.6
$ws.Range("I38").Value = 308.52173
$ws.Range("J38").Value = 34.57143
$ws.Range("K38").Value = 925.56519
$ws.Range("L38").Value = 103.71429
$ws.Range("M38").Value = -578.56519
$ws.Range("N38").Value = -797.71429
$ws.Range("H81").Value = 3049.5
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 3049.5
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 9148.5
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -11394.5
$ws.Range("H84").Value = 3049.5
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 3049.5
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 27445.5
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -38677.5
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H92").Value = 4933.6
$ws.Range("I92").Value = 3666.625
$ws.Range("J92").Value = 10001.5
$ws.Range("K92").Value = 10999.875
$ws.Range("L92").Value = 30004.5
$ws.Range("M92").Value = -9751.875
$ws.Range("N92").Value = -32500.5
$ws.Range("H132").Value = 4912.125
$ws.Range("J132").Value = 5969.8
$ws.Range("L132").Value = 53728.2
$ws.Range("N132").Value = -58788.2
$ws.Range("H137").Value = 4724
$ws.Range("I137").Value = 4724
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 14172
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -9072
$ws.Range("N137").ClearContents()
$ws.Range("H140").Value = 2590
$ws.Range("I140").Value = 1837
$ws.Range("K140").Value = 5511
$ws.Range("M140").Value = -331

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H44").Value = 4000
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H49").Value = 37499.5
$ws.Range("J49").Value = 37499.5
$ws.Range("L49").Value = 37499.5
$ws.Range("N49").Value = -37867.5
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5422.727
$ws.Range("I46").Value = 1162.5
$ws.Range("J46").Value = 7857.143
$ws.Range("K46").Value = 1162.5
$ws.Range("L46").Value = 7857.143
$ws.Range("M46").Value = -974.5
$ws.Range("N46").Value = -8233.143
$ws.Range("H100").Value = 6051
$ws.Range("I100").Value = 3097.4
$ws.Range("K100").Value = 3097.4
$ws.Range("M100").Value = -2556.4
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 7625.8
$ws.Range("I132").Value = 6560.263
$ws.Range("K132").Value = 19680.789
$ws.Range("M132").Value = -17150.789

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1004.4
$ws.Range("J100").Value = 916.3333
$ws.Range("L100").Value = 1832.6666
$ws.Range("N100").Value = -2914.6666
$ws.Range("H122").Value = 2791.9
$ws.Range("J122").Value = 4494
$ws.Range("L122").Value = 13482
$ws.Range("N122").Value = -18382
